# Fruta / hortaliza, semanal
# Insert two new weekly price rows (new rows 280 and 281) for
# "Vega Central Mapocho de Santiago - Arándano (blue)", pushing the
# previously-existing rows 280-334 down to 282-336.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 280..334 down by two, inserting two blank rows.
$ws.Rows("280:281").Insert()

# New row 280
$row280 = @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 45211, 13, "Fruta", 100101, "Berries", 100101001, "Arándano (blue)", "Sin especificar", "Primera", 320, 10000, 11000, 10375, "`$/bandeja 12 canastillos 125 gramos", "Provincia de Linares", 6917, 1.5)
for ($i = 0; $i -lt $row280.Length; $i++) {
    $ws.Cells.Item(280, $i + 1).Value = $row280[$i]
}

# New row 281
$row281 = @(9, "Vega Central Mapocho de Santiago", "Metropolitana", 45211, 13, "Fruta", 100101, "Berries", 100101001, "Arándano (blue)", "Sin especificar", "Primera", 160, 11000, 11000, 11000, "`$/bandeja 2 kilos", "Provincia de Curicó", 5500, 2)
for ($i = 0; $i -lt $row281.Length; $i++) {
    $ws.Cells.Item(281, $i + 1).Value = $row281[$i]
}
